$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = 'Volume 31   Number  2'
$ws.Range("C9").Value = 'Report Covering the Week  1/8/2024  Through  1/14/2024'
$ws.Range("D15").Value = 1
$ws.Range("F16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("N14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("F16").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100
$ws.Range("N14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1
$ws.Range("F16").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("K15").Value = -100
$ws.Range("N14").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 46.153846153846
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 12.5
$ws.Range("L16").Value = 125
$ws.Range("M16").Value = 28.571428571428
$ws.Range("N16").Value = -77.5
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -70
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -37.142857142857
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = -68
$ws.Range("L17").Value = -11.111111111111
$ws.Range("M17").Value = 33.333333333333
$ws.Range("N17").Value = -68
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = -66.666666666666
$ws.Range("N18").Value = -95.121951219512
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -51.020408163265
$ws.Range("I19").Value = 10
$ws.Range("J19").Value = 28
$ws.Range("L19").Value = -37.5
$ws.Range("M19").Value = -37.5
$ws.Range("N19").Value = -44.444444444444
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("F16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = 0
$ws.Range("N14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 116.666666666667
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 3
$ws.Range("F16").Copy()
$ws.Range("J20").PasteSpecial(-4122)
$ws.Range("K20").Value = 66.666666666666
$ws.Range("N14").Copy()
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -86.111111111111
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -52.941176470588
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -23.214285714285
$ws.Range("I21").Value = 34
$ws.Range("J21").Value = 71
$ws.Range("K21").Value = -52.112676056338
$ws.Range("L21").Value = -8.108108108108
$ws.Range("M21").Value = -2.857142857142
$ws.Range("N21").Value = -79.012345679012
$ws.Range("G22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = '***.*'
$ws.Range("C14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("L22").Value = -100
$ws.Range("N14").Copy()
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = -5.555555555555
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 12
$ws.Range("K23").Value = -41.666666666666
$ws.Range("L23").Value = 40
$ws.Range("M23").Value = 75
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 8.955223880597
$ws.Range("I24").Value = 30
$ws.Range("J24").Value = 33
$ws.Range("K24").Value = -9.090909090909
$ws.Range("L24").Value = -9.090909090909
$ws.Range("M24").Value = -37.5
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 55.263157894736
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 21
$ws.Range("K25").Value = 33.333333333333
$ws.Range("L25").Value = 100
$ws.Range("M25").Value = 64.705882352941
$ws.Range("D26").Value = 1
$ws.Range("F16").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100
$ws.Range("N14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("G26").Value = 1
$ws.Range("F16").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("H26").Value = -100
$ws.Range("N14").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("J26").Value = 1
$ws.Range("F16").Copy()
$ws.Range("J26").PasteSpecial(-4122)
$ws.Range("K26").Value = -100
$ws.Range("N14").Copy()
$ws.Range("K26").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -70
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -83.333333333333
$ws.Range("L27").Value = 0
$ws.Range("N14").Copy()
$ws.Range("L27").PasteSpecial(-4122)
$ws.Range("F28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("H28").Value = -100
$ws.Range("F29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("J36").Value = 6
$ws.Range("K36").Value = -45.454545454545
$ws.Range("L36").Value = -62.5
$ws.Range("M36").Value = -70
$ws.Range("N36").Value = -71.428571428571
$ws.Range("J43").Value = 1255
$ws.Range("K43").Value = -30.739514348785
$ws.Range("L43").Value = -54.627621113521
$ws.Range("M43").Value = -73.461619792768
$ws.Range("N43").Value = -79.580214773836
